# Updated sequence for bionwire wide trench design
#
# The "GoldWireBinder" handle set (columns B:H on row 8 of the Sequences,
# Names and Descriptions tabs) is re-ordered, and the Sequences row gets a
# shortened spacer ("TTTTG" dropped right after the "tt" linker) as part of
# the wide-trench BioNWire redesign.

$wb  = $excel.ActiveWorkbook
$seq = $wb.Worksheets.Item("Sequences")
$nam = $wb.Worksheets.Item("Names")
$des = $wb.Worksheets.Item("Descriptions")

# ---- Sequences (row 8, B:H) --------------------------------------------
$seq.Range("B8").Value = "CATCCTAATTCCGGTATTCTAAGAACGCTTCTGAATAATGGAttGGTTGATAAAAGCATGACAGGTTGATAATATAGAT"
$seq.Range("C8").Value = "GGCAAGGCATAGGTAAAGATTCAAAAGGCCGCCAGCCATTGCttGGTTGATAAAAGCATGACAGGTTGATAATATAGAT"
$seq.Range("D8").Value = "GAAACAATCGGCAAGAGACGCAGAAACAGCCGCACAGGCGGCttGGTTGATAAAAGCATGACAGGTTGATAATATAGAT"
$seq.Range("E8").Value = "TCAACATTAAATGGCGCATCGTAACCGTGCGGAAACCAGGCAttGGTTGATAAAAGCATGACAGGTTGATAATATAGAT"
$seq.Range("F8").Value = "TAGCAAGGCCGGCGTTTTCATCGGCATTTTCAGAGCCGCCACttGGTTGATAAAAGCATGACAGGTTGATAATATAGAT"
$seq.Range("G8").Value = "AACCTCCCGTTTTTGTTTAACGTCAAAAGATGGCAATTCATCttGGTTGATAAAAGCATGACAGGTTGATAATATAGAT"
$seq.Range("H8").Value = "CTTTACAGAGAAGCCCTTTTTAAGAAAACCAGAAGGAGCGGAttGGTTGATAAAAGCATGACAGGTTGATAATATAGAT"

# ---- Names (row 8, B:H) --------------------------------------------------
$nam.Range("B8").Value = "GoldWireBinder_h5_pos28"
$nam.Range("C8").Value = "GoldWireBinder_h5_pos11"
$nam.Range("D8").Value = "GoldWireBinder_h2_pos6"
$nam.Range("E8").Value = "GoldWireBinder_h2_pos8"
$nam.Range("F8").Value = "GoldWireBinder_h2_pos23"
$nam.Range("G8").Value = "GoldWireBinder_h5_pos27"
$nam.Range("H8").Value = "GoldWireBinder_h5_pos26"

# ---- Descriptions (row 8, B:H) -------------------------------------------
$des.Range("B8").Value = "Binding handle for the BioNWire gold nanowires at position 28 on side 5."
$des.Range("C8").Value = "Binding handle for the BioNWire gold nanowires at position 11 on side 5."
$des.Range("D8").Value = "Binding handle for the BioNWire gold nanowires at position 6 on side 2."
$des.Range("E8").Value = "Binding handle for the BioNWire gold nanowires at position 8 on side 2."
$des.Range("F8").Value = "Binding handle for the BioNWire gold nanowires at position 23 on side 2."
$des.Range("G8").Value = "Binding handle for the BioNWire gold nanowires at position 27 on side 5."
$des.Range("H8").Value = "Binding handle for the BioNWire gold nanowires at position 26 on side 5."

# ---- View state: per-sheet zoom + selection, matching final screen state --
# (Sequences keeps its existing zoom; Names is re-zoomed 150 -> 117;
#  Descriptions keeps its existing zoom but becomes the active tab.)
$seq.Activate()
$seq.Range("I8").Select()

$nam.Activate()
$excel.ActiveWindow.Zoom = 117
$nam.Range("B8:H8").Select()

$des.Activate()
$des.Range("D12").Select()
